$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: B2 "<them>" -> "<then>", C2 6 -> 7
$ws.Range("B2").Value = "<then>"
$ws.Range("C2").Value = 7

# Row 3: C3 5 -> 7
$ws.Range("C3").Value = 7

# Row 4: C4 5 -> 7
$ws.Range("C4").Value = 7

# Row 5: C5 12 -> 11
$ws.Range("C5").Value = 11

# Row 6: C6 6 -> 9
$ws.Range("C6").Value = 9

# Row 7: C7 2 -> 4
$ws.Range("C7").Value = 4

# Row 8: B8 "<november>" -> "<nomer>", C8 9 -> 11
$ws.Range("B8").Value = "<nomer>"
$ws.Range("C8").Value = 11

# Row 9: C9 6 -> 7
$ws.Range("C9").Value = 7

# Row 10: C10 8 -> 6
$ws.Range("C10").Value = 6

# Row 11: C11 9 -> 7
$ws.Range("C11").Value = 7

# Row 12: C12 7 -> 9
$ws.Range("C12").Value = 9

# Row 15: C15 6 -> 9
$ws.Range("C15").Value = 9

# Row 17: B17 "<would>" -> "<like>", C17 7 -> 5
$ws.Range("B17").Value = "<like>"
$ws.Range("C17").Value = 5

# Row 18: C18 6 -> 4
$ws.Range("C18").Value = 4
